$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings
# (e.g. "21.49") are not silently coerced into Excel numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.590.62'
$ws.Range('E2').Value = '  +3.05%  '
$ws.Range('D3').Value = '1.848.15'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  +3.09%  '
$ws.Range('D5').Value = '321.24'
$ws.Range('E5').Value = '  +3.53%  '
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('D7').Value = '0.4372'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('D8').Value = '0.3739'
$ws.Range('E8').Value = '  +2.13%  '
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('D10').Value = '0.8750'
$ws.Range('E10').Value = '  +1.64%  '
$ws.Range('D11').Value = '21.49'
$ws.Range('E11').Value = '  +3.21%  '
$ws.Range('D12').Value = '1.853.74'
$ws.Range('E12').Value = '  -3.36%  '
$ws.Range('D13').Value = '5.494'
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = '0.07155'
$ws.Range('E15').Value = '  +3.55%  '
$ws.Range('D16').Value = '82.80'
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('D18').Value = '0.000009021'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('E19').Value = '  +2.65%  '
$ws.Range('D20').Value = '15.41'
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('D21').Value = '27.581.67'
$ws.Range('E21').Value = '  +2.95%  '
$ws.Range('D22').Value = '5.242'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('D23').Value = '11.20'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').Value = '2.065.53'
$ws.Range('E24').Value = '  -3.55%  '
$ws.Range('E25').Value = '  +3.17%  '
$ws.Range('D26').Value = '1.931'
$ws.Range('E26').Value = '  +4.68%  '
$ws.Range('D27').Value = '18.71'
$ws.Range('E27').Value = '  +2.83%  '
$ws.Range('D28').Value = '5.267'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('D29').Value = '1.953'
$ws.Range('E29').Value = '  +2.62%  '
$ws.Range('D30').Value = '116.01'
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').Value = '0.09072'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = '1.207'
$ws.Range('E32').Value = '  +3.99%  '
$ws.Range('D33').Value = '0.7661'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').Value = '4.507'
$ws.Range('E34').Value = '  +2.25%  '
$ws.Range('D35').Value = '2.876'
$ws.Range('E35').Value = '  +3.94%  '
$ws.Range('D36').Value = '1.029'
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('D37').Value = '1.145'
$ws.Range('E37').Value = '  +2.35%  '
$ws.Range('E38').Value = '  +3.44%  '
$ws.Range('D39').Value = '0.05265'
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('D40').Value = '0.5163'
$ws.Range('E40').Value = '  +2.30%  '
$ws.Range('D41').Value = '2.802'
$ws.Range('E41').Value = '  +6.65%  '
$ws.Range('D42').Value = '0.1671'
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('B43').Value = 'PaxosStandard'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D43').Value = '1.120'
$ws.Range('E43').Value = '  +12.03%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '6.678'
$ws.Range('E44').Value = '  +2.90%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '8.555'
$ws.Range('E45').Value = '  +3.13%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '108.78'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '10.60'
$ws.Range('E47').Value = '  +3.15%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '1.710'
$ws.Range('E48').Value = '  +3.74%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = '0.4647'
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06372'
$ws.Range('E50').Value = '  +2.14%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.881'
$ws.Range('E51').Value = '  +4.79%  '
